$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.89786533337632157
$ws.Range("D2").Value = 0.13681624225895214
$ws.Range("E2").Value = 3.4391211244529551

# Row 3 (STR) values
$ws.Range("B3").Value = 0.28587858934978261
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 0.41173586391433881
$ws.Range("E3").Value = 1.9446659266568205

# Update selection to match new active range
$ws.Range("B1:E3").Select()
